$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Appended translation rows (smartInfo Create strings), rows 278-286 ---
$ws.Range("A278").Value = "Ist-Werte"
$ws.Range("B278").Value = "Actual-to-date"

$ws.Range("A279").Value = "Ist-Werte (Summe Portfolio)"
$ws.Range("B279").Value = "Actual-to-date (Sum OF Portfolio)"

$ws.Range("A280").Value = "Forecast (ETC)"
$ws.Range("B280").Value = "Forecast (ETC)"

$ws.Range("A281").Value = "Actual and Forecast (EAC)"
$ws.Range("B281").Value = "Actual and Forecast (EAC)"

$ws.Range("B282").Value = "Last ETC"
$ws.Range("A282").Value = "Letzter ETC"

$ws.Range("B283").Value = "Monthly Needs (Sum of Portfolio)"
$ws.Range("A283").Value = "monatl. Bedarfe (alle Projekte des Portfolios)"

$ws.Range("A284").Value = "Baseline (BAC)"
$ws.Range("B284").Value = "Baseline (BAC)"

$ws.Range("B285").Value = "Order-Value"
$ws.Range("A285").Value = "Auftragswert"

$ws.Range("B286").Value = "Invoices (Baseline)"
$ws.Range("A286").Value = "Rechnungen (Baseline)"

# --- More appended rows (repMessages now from VCSetting), rows 288-292 ---
$ws.Range("A288").Value = "Gesamt Kapazität"
$ws.Range("B288").Value = "Total Capacity"

$ws.Range("A289").Value = "interne Kapazität"
$ws.Range("B289").Value = "intern Capacity"

$ws.Range("A290").Value = "Summe Portfolio"
$ws.Range("B290").Value = "Sum of Portfolio"

$ws.Range("A291").Value = "Summe interner Mitarbeiter"
$ws.Range("B291").Value = "Sum of all intern employees"

$ws.Range("A292").Value = "Summe pro Monat"
$ws.Range("B292").Value = "Monthly Sum"

# row 287 inserted last, between the two existing blocks
$ws.Range("A287").Value = "Rechnungen (akt. Plan)"
$ws.Range("B287").Value = "Invoices (cur. Plan)"

# --- Style cleanup: a few cells revert from a transient explicit style back to the default (no) style ---
$ws.Range("A4").ClearFormats()
$ws.Range("B13").ClearFormats()
$ws.Range("B71").ClearFormats()

# --- Update the view/selection to match where the author ended up editing ---
[void]$ws.Range("B296").Select()
